$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Crops")
$win = $excel.ActiveWindow
$p = $win.Panes.Item(4)
Write-Host "ScrollRow:" $p.ScrollRow
Write-Host "ScrollColumn:" $p.ScrollColumn
$p.ScrollRow = 20
$p.ScrollColumn = 2
Write-Host "after ScrollRow:" $p.ScrollRow
Write-Host "after ScrollColumn:" $p.ScrollColumn
